{"js": "const newValues = [\"37+29=66\", \"51-7=44\", \"29+68=97\", \"70-55=15\", \"61-25=36\", \"72-64=8\", \"91-73=18\", \"84+8=92\", \"52-38=14\", \"35+46=81\", \"17+36=53\", \"35+47=82\", \"3+18=21\", \"93-28=65\", \"57+28=85\", \"43-19=24\", \"3+69=72\", \"41-17=24\", \"90-19=71\", \"81-75=6\", \"90-65=25\", \"17+36=53\", \"48+16=64\", \"76+9=85\", \"42+9=51\", \"60-39=21\", \"51-23=28\", \"28+45=73\", \"93-47=46\", \"28+19=47\", \"9+85=94\", \"30-17=13\", \"60-37=23\", \"85-29=56\", \"82-24=58\", \"58+27=85\", \"43+39=82\", \"90-87=3\", \"57+18=75\", \"91-89=2\", \"37+19=56\", \"68-29=39\", \"28+46=74\", \"71-22=49\", \"94-66=28\", \"7+59=66\", \"29+13=42\", \"93-74=19\", \"96-28=68\", \"16+39=55\", \"37-8=29\", \"92-45=47\", \"81-52=29\", \"50-49=1\", \"29+18=47\", \"84-9=75\", \"37+47=84\", \"92-43=49\", \"53-34=19\", \"24+18=42\", \"43+28=71\", \"40-11=29\", \"66-59=7\", \"71-6=65\", \"48+45=93\", \"2+59=61\", \"86+7=93\", \"58+29=87\", \"89+4=93\", \"40-37=3\", \"76+8=84\", \"56-29=27\", \"15+17=32\", \"31-6=25\", \"36+17=53\", \"63+28=91\", \"78+6=84\", \"69+14=83\", \"33-8=25\", \"3+88=91\", \"88+8=96\", \"19+56=75\", \"33+58=91\", \"37+54=91\", \"22-14=8\", \"81-49=32\", \"73-44=29\", \"82-3=79\", \"48+47=95\", \"90-16=74\", \"65+9=74\", \"6+19=25\", \"3+49=52\", \"88-69=19\", \"76-39=37\", \"45+8=53\", \"16+65=81\", \"3+79=82\", \"70-29=41\", \"36+5=41\"];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst cols = 5;\nlet idx = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < cols; c++) {\n    const cell = table.getCell(r, c);\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n    const p = paragraphs.items[0];\n    p.insertText(newValues[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n  await context.sync();\n}\n", "ps1": "$newValues = @(\"37+29=66\",\"51-7=44\",\"29+68=97\",\"70-55=15\",\"61-25=36\",\"72-64=8\",\"91-73=18\",\"84+8=92\",\"52-38=14\",\"35+46=81\",\"17+36=53\",\"35+47=82\",\"3+18=21\",\"93-28=65\",\"57+28=85\",\"43-19=24\",\"3+69=72\",\"41-17=24\",\"90-19=71\",\"81-75=6\",\"90-65=25\",\"17+36=53\",\"48+16=64\",\"76+9=85\",\"42+9=51\",\"60-39=21\",\"51-23=28\",\"28+45=73\",\"93-47=46\",\"28+19=47\",\"9+85=94\",\"30-17=13\",\"60-37=23\",\"85-29=56\",\"82-24=58\",\"58+27=85\",\"43+39=82\",\"90-87=3\",\"57+18=75\",\"91-89=2\",\"37+19=56\",\"68-29=39\",\"28+46=74\",\"71-22=49\",\"94-66=28\",\"7+59=66\",\"29+13=42\",\"93-74=19\",\"96-28=68\",\"16+39=55\",\"37-8=29\",\"92-45=47\",\"81-52=29\",\"50-49=1\",\"29+18=47\",\"84-9=75\",\"37+47=84\",\"92-43=49\",\"53-34=19\",\"24+18=42\",\"43+28=71\",\"40-11=29\",\"66-59=7\",\"71-6=65\",\"48+45=93\",\"2+59=61\",\"86+7=93\",\"58+29=87\",\"89+4=93\",\"40-37=3\",\"76+8=84\",\"56-29=27\",\"15+17=32\",\"31-6=25\",\"36+17=53\",\"63+28=91\",\"78+6=84\",\"69+14=83\",\"33-8=25\",\"3+88=91\",\"88+8=96\",\"19+56=75\",\"33+58=91\",\"37+54=91\",\"22-14=8\",\"81-49=32\",\"73-44=29\",\"82-3=79\",\"48+47=95\",\"90-16=74\",\"65+9=74\",\"6+19=25\",\"3+49=52\",\"88-69=19\",\"76-39=37\",\"45+8=53\",\"16+65=81\",\"3+79=82\",\"70-29=41\",\"36+5=41\")\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$cols = 5\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$idx]\n    $idx++\n  }\n}\n"}
